$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(62).Insert()

$ws.Cells.Item(62, 1).Value = 9
$ws.Cells.Item(62, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = 45259
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100101
$ws.Cells.Item(62, 8).Value = "Berries"
$ws.Cells.Item(62, 9).Value = 100101004
$ws.Cells.Item(62, 10).Value = "Frambuesa"
$ws.Cells.Item(62, 11).Value = "Sin especificar"
$ws.Cells.Item(62, 12).Value = "Especial"
$ws.Cells.Item(62, 13).Value = 300
$ws.Cells.Item(62, 14).Value = 10000
$ws.Cells.Item(62, 15).Value = 10000
$ws.Cells.Item(62, 16).Value = 10000
$ws.Cells.Item(62, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(62, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(62, 19).Value = 5000
$ws.Cells.Item(62, 20).Value = 2
